$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2968.0205
$ws.Range("J138").Value = 3040.5715
$ws.Range("L138").Value = 9121.7145
$ws.Range("N138").Value = -19401.7145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 390823.6
$ws.Range("I32").Value = 521457.62
$ws.Range("J32").Value = 12672.421
$ws.Range("K32").Value = 521457.62
$ws.Range("L32").Value = 12672.421
$ws.Range("M32").Value = -521170.62
$ws.Range("N32").Value = -13246.421
$ws.Range("H61").Value = 7754032
$ws.Range("I61").Value = 12347115
$ws.Range("K61").Value = 12347115
$ws.Range("M61").Value = -12346903
$ws.Range("H74").Value = 715
$ws.Range("I74").Value = 569.1875
$ws.Range("J74").Value = 816.43475
$ws.Range("K74").Value = 569.1875
$ws.Range("L74").Value = 816.43475
$ws.Range("M74").Value = 304.8125
$ws.Range("N74").Value = -2564.43475
$ws.Range("H77").Value = 715
$ws.Range("I77").Value = 569.1875
$ws.Range("J77").Value = 816.43475
$ws.Range("K77").Value = 2845.9375
$ws.Range("L77").Value = 4082.17375
$ws.Range("M77").Value = 1522.0625
$ws.Range("N77").Value = -12818.17375
$ws.Range("H121").Value = 47490
$ws.Range("J121").Value = 47490
$ws.Range("L121").Value = 47490
$ws.Range("N121").Value = -50984
$ws.Range("H136").Value = 7754032
$ws.Range("I136").Value = 12347115
$ws.Range("K136").Value = 37041345
$ws.Range("M136").Value = -37038795

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 96390
$ws.Range("J55").Value = 96390
$ws.Range("L55").Value = 96390
$ws.Range("N55").Value = -96936
$ws.Range("H58").Value = 19800
$ws.Range("J58").Value = 19800
$ws.Range("L58").Value = 19800
$ws.Range("N58").Value = -20388
$ws.Range("H123").Value = 99780
$ws.Range("J123").Value = 99780
$ws.Range("L123").Value = 99780
$ws.Range("N123").Value = -109580
$ws.Range("H134").Value = 3072.7334
$ws.Range("I134").Value = 3019.1428
$ws.Range("K134").Value = 9057.428400000001
$ws.Range("M134").Value = -6522.428400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1296.4889
$ws.Range("I31").Value = 1038.725
$ws.Range("J31").Value = 3358.6
$ws.Range("K31").Value = 1038.725
$ws.Range("L31").Value = 3358.6
$ws.Range("M31").Value = -743.7249999999999
$ws.Range("N31").Value = -3948.6
$ws.Range("H34").Value = 1296.4889
$ws.Range("I34").Value = 1038.725
$ws.Range("J34").Value = 3358.6
$ws.Range("K34").Value = 1038.725
$ws.Range("L34").Value = 3358.6
$ws.Range("M34").Value = -836.7249999999999
$ws.Range("N34").Value = -3762.6
$ws.Range("H58").Value = 3273.0588
$ws.Range("I58").Value = 3376.5454
$ws.Range("J58").Value = 3083.3333
$ws.Range("K58").Value = 3376.5454
$ws.Range("L58").Value = 3083.3333
$ws.Range("M58").Value = -3173.5454
$ws.Range("N58").Value = -3489.3333
$ws.Range("H74").Value = 28285.428
$ws.Range("J74").Value = 28285.428
$ws.Range("L74").Value = 28285.428
$ws.Range("N74").Value = -30033.428
$ws.Range("H77").Value = 28285.428
$ws.Range("J77").Value = 28285.428
$ws.Range("L77").Value = 84856.284
$ws.Range("N77").Value = -93592.284
$ws.Range("H132").Value = 9261477
$ws.Range("I132").Value = 1910.6
$ws.Range("J132").Value = 20835936
$ws.Range("K132").Value = 5731.799999999999
$ws.Range("L132").Value = 62507808
$ws.Range("M132").Value = -3201.799999999999
$ws.Range("N132").Value = -62512868
$ws.Range("H134").Value = 967.5833
$ws.Range("I134").Value = 768
$ws.Range("J134").Value = 1566.3334
$ws.Range("K134").Value = 2304
$ws.Range("L134").Value = 4699.0002
$ws.Range("M134").Value = 231
$ws.Range("N134").Value = -9769.0002
$ws.Range("H136").Value = 3273.0588
$ws.Range("I136").Value = 3376.5454
$ws.Range("J136").Value = 3083.3333
$ws.Range("K136").Value = 10129.6362
$ws.Range("L136").Value = 9249.999899999999
$ws.Range("M136").Value = -7579.636200000001
$ws.Range("N136").Value = -14349.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 285.7143
$ws.Range("I12").Value = 200
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 900
$ws.Range("M12").Value = -427
$ws.Range("N12").Value = -1246
$ws.Range("H68").Value = 778.62
$ws.Range("I68").Value = 695
$ws.Range("J68").Value = 780.32654
$ws.Range("K68").Value = 2085
$ws.Range("L68").Value = 2340.97962
$ws.Range("M68").Value = -1274
$ws.Range("N68").Value = -3962.97962
$ws.Range("H71").Value = 778.62
$ws.Range("I71").Value = 695
$ws.Range("J71").Value = 780.32654
$ws.Range("K71").Value = 6255
$ws.Range("L71").Value = 7022.93886
$ws.Range("M71").Value = -2199
$ws.Range("N71").Value = -15134.93886
$ws.Range("H107").Value = 1468.8823
$ws.Range("I107").Value = 268.4643
$ws.Range("J107").Value = 2930.261
$ws.Range("K107").Value = 805.3928999999999
$ws.Range("L107").Value = 8790.782999999999
$ws.Range("M107").Value = 1114.6071
$ws.Range("N107").Value = -12630.783
$ws.Range("H122").Value = 5627.45
$ws.Range("I122").Value = 347.625
$ws.Range("J122").Value = 9147.333000000001
$ws.Range("K122").Value = 3128.625
$ws.Range("L122").Value = 82325.997
$ws.Range("M122").Value = -678.625
$ws.Range("N122").Value = -87225.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H116").Value = 60000
$ws.Range("J116").Value = 60000
$ws.Range("L116").Value = 60000
$ws.Range("N116").Value = -69178
$ws.Range("H122").Value = 1642
$ws.Range("I122").Value = 1410.1428
$ws.Range("K122").Value = 4230.428400000001
$ws.Range("M122").Value = -1780.428400000001
$ws.Range("H132").Value = 2713.7144
$ws.Range("I132").Value = 2390.6956
$ws.Range("K132").Value = 7172.0868
$ws.Range("M132").Value = -4642.0868

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1455.2858
$ws.Range("I136").Value = 1452.875
$ws.Range("K136").Value = 4358.625
$ws.Range("M136").Value = -1808.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5590.5713
$ws.Range("I81").Value = 4855.8335
$ws.Range("J81").Value = 9999
$ws.Range("K81").Value = 9711.666999999999
$ws.Range("L81").Value = 19998
$ws.Range("M81").Value = -8650.666999999999
$ws.Range("N81").Value = -22120
$ws.Range("H84").Value = 5590.5713
$ws.Range("I84").Value = 4855.8335
$ws.Range("J84").Value = 9999
$ws.Range("K84").Value = 48558.335
$ws.Range("L84").Value = 99990
$ws.Range("M84").Value = -43254.335
$ws.Range("N84").Value = -110598
$ws.Range("H132").Value = 7677376
$ws.Range("I132").Value = 1763.6296
$ws.Range("J132").Value = 26517516
$ws.Range("K132").Value = 5290.8888
$ws.Range("L132").Value = 79552548
$ws.Range("M132").Value = -2760.8888
$ws.Range("N132").Value = -79557608

Write-Host "Applied all cell updates"